$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.736.73'
$ws.Range("E2").Value = '  -0.23%  '

# Row 3
$ws.Range("D3").Value = '1.630.76'
$ws.Range("E3").Value = '  -0.37%  '

# Row 4
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").Value = '''214.53'
$ws.Range("E5").Value = '  -0.45%  '

# Row 6
$ws.Range("E6").Value = '  -0.73%  '

# Row 7
$ws.Range("E7").Value = '  -0.10%  '

# Row 8
$ws.Range("D8").Value = '''0.255'
$ws.Range("E8").Value = '  -1.22%  '

# Row 9
$ws.Range("D9").Value = '''0.0633'
$ws.Range("E9").Value = '  -1.33%  '

# Row 10
$ws.Range("D10").Value = '''19.54'
$ws.Range("E10").Value = '  -1.63%  '

# Row 11
$ws.Range("E11").Value = '  +1.57%  '

# Row 12
$ws.Range("D12").Value = '''4.26'

# Row 13
$ws.Range("D13").Value = '1.855.38'
$ws.Range("E13").Value = '  -0.28%  '

# Row 14
$ws.Range("D14").Value = '1.611.66'
$ws.Range("E14").Value = '  -1.62%  '

# Row 15
$ws.Range("D15").Value = '''0.556'
$ws.Range("E15").Value = '  -0.17%  '

# Row 16
$ws.Range("D16").Value = '0.0₃0761'
$ws.Range("E16").Value = '  -2.11%  '

# Row 17
$ws.Range("D17").Value = '''62.93'
$ws.Range("E17").Value = '  -0.51%  '

# Row 18
$ws.Range("D18").Value = '25.728.27'
$ws.Range("E18").Value = '  -0.30%  '

# Row 19
$ws.Range("E19").Value = '  -0.07%  '

# Row 20
$ws.Range("D20").Value = '''4.44'
$ws.Range("E20").Value = '  -0.12%  '

# Row 21
$ws.Range("D21").Value = '''192.09'
$ws.Range("E21").Value = '  -1.16%  '

# Row 22
$ws.Range("D22").Value = '''9.93'
$ws.Range("E22").Value = '  -0.50%  '

# Row 23
$ws.Range("D23").Value = '''6.25'
$ws.Range("E23").Value = '  +1.41%  '

# Row 24
$ws.Range("E24").Value = '  -0.11%  '

# Row 25
$ws.Range("D25").Value = '''1.82'
$ws.Range("E25").Value = '  +3.14%  '

# Row 26
$ws.Range("D26").Value = '''142.29'
$ws.Range("E26").Value = '  +1.61%  '

# Row 27
$ws.Range("E27").Value = '  +1.83%  '

# Row 28
$ws.Range("D28").Value = '''6.87'
$ws.Range("E28").Value = '  +0.05%  '

# Row 29
$ws.Range("D29").Value = '''15.47'
$ws.Range("E29").Value = '  -1.11%  '

# Row 30
$ws.Range("E30").Value = '  -0.52%  '

# Row 31
$ws.Range("D31").Value = '''0.0490'
$ws.Range("E31").Value = '  -0.45%  '

# Row 32
$ws.Range("E32").Value = '  -0.78%  '

# Row 33
$ws.Range("D33").Value = '''3.23'
$ws.Range("E33").Value = '  -1.26%  '

# Row 34
$ws.Range("E34").Value = '  -0.95%  '

# Row 35
$ws.Range("E35").Value = '  +0.51%  '

# Row 36
$ws.Range("D36").Value = '''0.905'
$ws.Range("E36").Value = '  +0.56%  '

# Row 37
$ws.Range("D37").Value = '1.140.25'
$ws.Range("E37").Value = '  +2.79%  '

# Row 38
$ws.Range("E38").Value = '  -2.58%  '

# Row 39
$ws.Range("D39").Value = '''0.543'
$ws.Range("E39").Value = '  -1.88%  '

# Row 40
$ws.Range("E40").Value = '  -0.98%  '

# Row 41
$ws.Range("D41").Value = '''2.55'
$ws.Range("E41").Value = '  -0.18%  '

# Row 42
$ws.Range("D42").Value = '''0.999'
$ws.Range("E42").Value = '  -0.13%  '

# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.53'
$ws.Range("E43").Value = '  -0.67%  '

# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''100.46'
$ws.Range("E44").Value = '  +1.23%  '

# Row 45
$ws.Range("D45").Value = '''0.803'
$ws.Range("E45").Value = '  -0.19%  '

# Row 46
$ws.Range("D46").Value = '1.764.76'
$ws.Range("E46").Value = '  -0.03%  '

# Row 47
$ws.Range("E47").Value = '  +0.51%  '

# Row 48
$ws.Range("D48").Value = '''55.27'
$ws.Range("E48").Value = '  -0.10%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.0507'
$ws.Range("E49").Value = '  +0.68%  '

# Row 50
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '''0.418'
$ws.Range("E50").Value = '  +0.12%  '

# Row 51
$ws.Range("D51").Value = '''1.44'
$ws.Range("E51").Value = '  +4.32%  '
